$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.914.47"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "'2.322.83"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'302.61"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'96.55"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'34.65"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "'19.12"
$ws.Range("E11").Value = "  +6.53%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "'6.79"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "'2.687.84"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'2.325.04"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "'42.852.67"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "'12.21"
$ws.Range("E19").Value = "  -5.08%  "
$ws.Range("D20").Value = "'6.17"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "'68.09"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("E23").Value = "  +6.23%  "
$ws.Range("D24").Value = "'236.46"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").Value = "'24.42"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "'166.25"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'9.15"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "'32.44"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'5.02"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "'17.83"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "'4.49"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'0.0702"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "'0.1000"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "'20.90"
$ws.Range("E42").Value = "  +14.26%  "
$ws.Range("D43").Value = "'1.936.11"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "'10.29"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("D47").Value = "'2.77"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").Value = "'2.555.05"
$ws.Range("D49").Value = "'53.55"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'2.79"
$ws.Range("E50").Value = "  -4.57%  "
$ws.Range("D51").Value = "'72.20"
$ws.Range("E51").Value = "  +2.58%  "
